$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.757.30'
$ws.Range("E2").Value = '  -6.70%  '

$ws.Range("D3").Value = '1.697.24'
$ws.Range("E3").Value = '  -5.86%  '

$ws.Range("D4").Value = "'" + '1.006'
$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").Value = "'" + '220.16'
$ws.Range("E5").Value = '  -4.88%  '

$ws.Range("D6").Value = "'" + '0.5079'
$ws.Range("E6").Value = '  -14.68%  '

$ws.Range("D7").Value = "'" + '1.006'
$ws.Range("E7").Value = '  +0.33%  '

$ws.Range("D8").Value = "'" + '0.2601'
$ws.Range("E8").Value = '  -6.08%  '

$ws.Range("D9").Value = "'" + '21.94'
$ws.Range("E9").Value = '  -6.00%  '

$ws.Range("D10").Value = "'" + '0.06147'
$ws.Range("E10").Value = '  -9.83%  '

$ws.Range("D11").Value = "'" + '0.07356'
$ws.Range("E11").Value = '  -2.08%  '

$ws.Range("D12").Value = '1.669.69'
$ws.Range("E12").Value = '  -7.47%  '

$ws.Range("D13").Value = "'" + '4.462'
$ws.Range("E13").Value = '  -4.94%  '

$ws.Range("D14").Value = "'" + '0.5761'
$ws.Range("E14").Value = '  -8.28%  '

$ws.Range("D15").Value = '1.929.68'
$ws.Range("E15").Value = '  -5.82%  '

$ws.Range("D16").Value = "'" + '0.000008137'
$ws.Range("E16").Value = '  -11.63%  '

$ws.Range("D17").Value = "'" + '65.69'
$ws.Range("E17").Value = '  -12.95%  '

$ws.Range("D18").Value = '26.798.03'
$ws.Range("E18").Value = '  -6.23%  '

$ws.Range("D19").Value = "'" + '5.032'
$ws.Range("E19").Value = '  -8.07%  '

$ws.Range("D20").Value = "'" + '1.006'
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("D21").Value = "'" + '10.74'
$ws.Range("E21").Value = '  -5.69%  '

$ws.Range("D22").Value = "'" + '183.55'
$ws.Range("E22").Value = '  -12.23%  '

$ws.Range("D23").Value = "'" + '6.230'
$ws.Range("E23").Value = '  -8.47%  '

$ws.Range("D24").Value = "'" + '1.007'
$ws.Range("E24").Value = '  +0.35%  '

$ws.Range("D25").Value = "'" + '145.21'
$ws.Range("E25").Value = '  -6.04%  '

$ws.Range("D26").Value = "'" + '7.617'
$ws.Range("E26").Value = '  -2.99%  '

$ws.Range("D27").Value = "'" + '0.1146'
$ws.Range("E27").Value = '  -9.95%  '

$ws.Range("D28").Value = "'" + '15.25'
$ws.Range("E28").Value = '  -6.94%  '

$ws.Range("D29").Value = "'" + '1.322'
$ws.Range("E29").Value = '  -8.63%  '

$ws.Range("D30").Value = "'" + '0.05669'
$ws.Range("E30").Value = '  -9.99%  '

$ws.Range("E31").Value = '  -5.13%  '

$ws.Range("D32").Value = "'" + '3.476'
$ws.Range("E32").Value = '  -7.14%  '

$ws.Range("D33").Value = "'" + '3.425'
$ws.Range("E33").Value = '  -7.99%  '

$ws.Range("D34").Value = "'" + '1.677'
$ws.Range("E34").Value = '  -2.40%  '

$ws.Range("D35").Value = "'" + '1.005'
$ws.Range("E35").Value = '  -4.26%  '

$ws.Range("D36").Value = "'" + '2.415'
$ws.Range("E36").Value = '  -3.62%  '

$ws.Range("D37").Value = "'" + '0.5936'
$ws.Range("E37").Value = '  -6.62%  '

$ws.Range("D38").Value = "'" + '2.646'
$ws.Range("E38").Value = '  -2.76%  '

$ws.Range("D39").Value = "'" + '0.01596'
$ws.Range("E39").Value = '  -6.17%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'" + '5.930'
$ws.Range("E40").Value = '  -7.42%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.073.65'
$ws.Range("E41").Value = '  -5.63%  '

$ws.Range("D42").Value = "'" + '0.8555'
$ws.Range("E42").Value = '  -0.64%  '

$ws.Range("D43").Value = "'" + '1.004'
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("D44").Value = "'" + '98.04'
$ws.Range("E44").Value = '  -2.83%  '

$ws.Range("D45").Value = '1.841.49'
$ws.Range("E45").Value = '  -6.07%  '

$ws.Range("D46").Value = "'" + '56.47'
$ws.Range("E46").Value = '  -6.79%  '

$ws.Range("E47").Value = '  -5.90%  '

$ws.Range("D48").Value = "'" + '1.003'
$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("D49").Value = "'" + '8.027'
$ws.Range("E49").Value = '  -3.32%  '

$ws.Range("D50").Value = "'" + '0.4343'
$ws.Range("E50").Value = '  -3.49%  '

$ws.Range("D51").Value = "'" + '0.05215'
$ws.Range("E51").Value = '  -4.27%  '
